$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for new columns I (I0) and J (IF), matching the
# existing header style (bold font, thin border, centered) used by H1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Populate data values for I2:J69
$values = @(
    @(2, 8, 9),
    @(3, 8, 8),
    @(4, 1, 1),
    @(5, 7, 7),
    @(6, 7, 8),
    @(7, 8, 8),
    @(8, 8, 8),
    @(9, 8, 8),
    @(10, 7, 8),
    @(11, 8, 8),
    @(12, 7, 8),
    @(13, 6, 7),
    @(14, 6, 7),
    @(15, 3, 4),
    @(16, 8, 8),
    @(17, 6, 6),
    @(18, 8, 9),
    @(19, 8, 9),
    @(20, 7, 8),
    @(21, 7, 8),
    @(22, 7, 7),
    @(23, 8, 8),
    @(24, 8, 8),
    @(25, 7, 7),
    @(26, 7, 7),
    @(27, 6, 7),
    @(28, 6, 6),
    @(29, 6, 6),
    @(30, 8, 8),
    @(31, 8, 8),
    @(32, 7, 7),
    @(33, 9, 9),
    @(34, 9, 9),
    @(35, 6, 6),
    @(36, 6, 6),
    @(37, 7, 7),
    @(38, 8, 8),
    @(39, 8, 8),
    @(40, 10, 11),
    @(41, 6, 7),
    @(42, 5, 6),
    @(43, 9, 9),
    @(44, 8, 8),
    @(45, 9, 9),
    @(46, 9, 9),
    @(47, 9, 9),
    @(48, 5, 6),
    @(49, 7, 7),
    @(50, 8, 8),
    @(51, 9, 9),
    @(52, 9, 9),
    @(53, 6, 6),
    @(54, 7, 7),
    @(55, 7, 8),
    @(56, 9, 9),
    @(57, 7, 7),
    @(58, 7, 8),
    @(59, 5, 5),
    @(60, 5, 5),
    @(61, 8, 8),
    @(62, 5, 6),
    @(63, 6, 6),
    @(64, 6, 7),
    @(65, 9, 9),
    @(66, 5, 6),
    @(67, 7, 7),
    @(68, 9, 9),
    @(69, 2, 2)
)

foreach ($row in $values) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
